$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$av = $excel.ActiveWindow
Write-Host ("Panes.Count=" + $av.Panes.Count)
$p1 = $av.Panes.Item(1)
Write-Host ("Pane1 ScrollRow=" + $p1.ScrollRow + " ScrollColumn=" + $p1.ScrollColumn)
$p2 = $av.Panes.Item(2)
Write-Host ("Pane2 ScrollRow=" + $p2.ScrollRow + " ScrollColumn=" + $p2.ScrollColumn)
